$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 15.84219932556152
$ws.Range("D2").Value = 177

$ws.Range("C3").Value = 15.38395881652832
$ws.Range("D3").Value = 176

$ws.Range("C4").Value = 21.89183235168457
$ws.Range("D4").Value = 176

$ws.Range("C5").Value = 16.91484451293945
$ws.Range("D5").Value = 176

$ws.Range("C6").Value = 16.65210723876953
$ws.Range("D6").Value = 177
